# Adds the "Metodos Norma 1, 2, inf en matriz" and
# "Metodos Norma 1, 2, inf en vector" rows to the Metricas sheet
# (practica2 TDAs - metricas_sel.xlsx), filling in the two previously
# empty data rows (8 and 9) of the task table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 8: Metodos Norma 1, 2, inf en matriz ---
$ws.Range("A8").Value = "Metodos Norma 1, 2, inf en matriz"
$ws.Range("B8").Value = 20
$ws.Range("C8").Value = 24
$ws.Range("D8").Value = 0.017361111111111112
$ws.Range("E8").Value = 0.7125
$ws.Range("F8").Value = 0.7256944444444445
$ws.Range("H8").Value = 1
$ws.Range("I8").Value = 0.0006944444444444445

# --- Row 9: Metodos Norma 1, 2, inf en vector ---
$ws.Range("A9").Value = "Metodos Norma 1, 2, inf en vector"
$ws.Range("B9").Value = 20
$ws.Range("D9").Value = 0.006944444444444444
$ws.Range("E9").Value = 0.7284722222222223
$ws.Range("F9").Value = 0
$ws.Range("H9").Value = 0
$ws.Range("I9").Value = 0

$ws.Range("I9").Select()
